$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A1").Value = 0.258999228477478
$ws.Range("B1").Value = 2.093186616897583
$ws.Range("C1").Value = 3.549638748168945
$ws.Range("D1").Value = 3.61506175994873
$ws.Range("E1").Value = 0.8363673090934753
